# STS script update: replace the sample/random placeholder data with the
# new "cseregy-test" rows, add a header row, and add a new (empty,
# auto-styled "Hyperlink") column C header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook used to hold 6 data rows in columns A:B; the new layout only
# needs 5 rows (2 header rows + 3 data rows), so drop the now-unused last row.
$ws.Rows(6).Delete()

# --- Row 1: column headers -------------------------------------------------
$ws.Range("A1").Value = "column_0"
$ws.Range("B1").Value = "column_1"

# --- Row 2: field-name headers ---------------------------------------------
$ws.Range("A2").Value = "resource_name"
$ws.Range("B2").Value = "account_number"

# --- Row 3-5: actual data ---------------------------------------------------
$ws.Range("A3").Value = "cseregy-test-1"
$ws.Range("B3").Value = 341254748130

$ws.Range("A4").Value = "cseregy-test-2"
$ws.Range("B4").Value = 341254748130

$ws.Range("A5").Value = "5945869c-1517-4c88-b5b2-b8d3fd35b2f4"
$ws.Range("B5").Value = 341254748130

# --- C1: a leftover hyperlink-formatted (but empty) cell -------------------
# A hyperlink was inserted and removed here, which is why Excel keeps the
# auto-generated "Hyperlink" cell style applied even though the cell itself
# ends up empty.
$ws.Hyperlinks.Add($ws.Range("C1"), "http://example.com", "", "", "link")
$ws.Hyperlinks(1).Delete()
$ws.Range("C1").ClearContents()

# --- Column widths (resized to fit the new, wider content) -----------------
# Column B's target width is essentially unchanged from the original, so it
# is left alone; columns A and C need to widen noticeably for the new values.
$ws.Columns("A").ColumnWidth = 26.7109375
$ws.Columns("C").ColumnWidth = 36.85546875

# --- Selection / active cell -------------------------------------------------
$null = $ws.Range("C4").Select()
